# Update gh-pages to output generated at 456a3b4
# This script updates the "想去人数" (F column) counts on the 4 worksheets
# to the freshly-generated values, matching the target XML diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2657
$ws1.Range("F3").Value = 576
$ws1.Range("F7").Value = 488
$ws1.Range("F8").Value = 1228
$ws1.Range("F9").Value = 572
$ws1.Range("F10").Value = 311
$ws1.Range("F11").Value = 5
$ws1.Range("F12").Value = 129
$ws1.Range("F13").Value = 361
$ws1.Range("F14").Value = 5758
$ws1.Range("F15").Value = 85
$ws1.Range("F16").Value = 1787
$ws1.Range("F17").Value = 4194
$ws1.Range("F18").Value = 437
$ws1.Range("F19").Value = 239
$ws1.Range("F20").Value = 303
$ws1.Range("F21").Value = 4896
$ws1.Range("F22").Value = 6268
$ws1.Range("F25").Value = 697
$ws1.Range("F26").Value = 3788
$ws1.Range("F27").Value = 502
$ws1.Range("F29").Value = 197
$ws1.Range("F32").Value = 1419
$ws1.Range("F34").Value = 569
$ws1.Range("F35").Value = 1608
$ws1.Range("F36").Value = 202
$ws1.Range("F37").Value = 1731
$ws1.Range("F38").Value = 202
$ws1.Range("F39").Value = 1146
$ws1.Range("F40").Value = 1336
$ws1.Range("F41").Value = 636
$ws1.Range("F42").Value = 97
$ws1.Range("F43").Value = 3435
$ws1.Range("F45").Value = 292
$ws1.Range("F46").Value = 417
$ws1.Range("F48").Value = 19
$ws1.Range("F49").Value = 3894

# Sheet "演出" (sheetId 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 1211
$ws2.Range("F16").Value = 13
$ws2.Range("F24").Value = 46

# Sheet "本地生活" (sheetId 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 3953

# Sheet "全部类型" (sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3953
$ws4.Range("F4").Value = 576
$ws4.Range("F7").Value = 1211
$ws4.Range("F11").Value = 488
$ws4.Range("F13").Value = 1228
$ws4.Range("F14").Value = 572
$ws4.Range("F15").Value = 311
$ws4.Range("F16").Value = 129
$ws4.Range("F17").Value = 361
$ws4.Range("F18").Value = 1787
$ws4.Range("F19").Value = 4194
$ws4.Range("F20").Value = 4896
$ws4.Range("F23").Value = 697
$ws4.Range("F24").Value = 3788
$ws4.Range("F25").Value = 502
$ws4.Range("F27").Value = 197
$ws4.Range("F29").Value = 1419
$ws4.Range("F31").Value = 569
$ws4.Range("F32").Value = 1608
$ws4.Range("F33").Value = 202
$ws4.Range("F34").Value = 1731
$ws4.Range("F37").Value = 636
$ws4.Range("F39").Value = 97
$ws4.Range("F41").Value = 3435
$ws4.Range("F44").Value = 292
$ws4.Range("F45").Value = 417
$ws4.Range("F48").Value = 3894
